$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "forgotten password" row's leftover placeholder values (A6/B6
# previously held a stray single-space string); row 6 now only carries the
# expected-result columns (C6/D6).
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()

# Move the saved selection to D10 (matches the resaved workbook's cursor
# position).
$ws.Range("D10").Select()
